$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "62.916.86"
$ws.Range("E2").Value = "  -1.76%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.679.13"
$ws.Range("E3").Value = "  -2.15%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "552.46"
$ws.Range("E5").Value = "  -3.11%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "158.33"
$ws.Range("E6").Value = "  -1.00%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.593"
$ws.Range("E8").Value = "  -0.69%  "
$ws.Range("E9").Value = "  -3.27%  "
$ws.Range("E10").Value = "  -1.46%  "
$ws.Range("E11").Value = "  -4.17%  "
$ws.Range("E12").Value = "  -7.07%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.151.69"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "26.25"
$ws.Range("E14").Value = "  -2.02%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "62.781.23"
$ws.Range("E15").Value = "  -1.32%  "
$ws.Range("E16").Value = "  -2.47%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.677.87"
$ws.Range("E17").Value = "  -2.30%  "
$ws.Range("E18").Value = "  -1.88%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.61"
$ws.Range("E19").Value = "  -4.21%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "345.31"
$ws.Range("E20").Value = "  -2.53%  "
$ws.Range("E21").Value = "  -4.65%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.00"
$ws.Range("E22").Value = "  +0.10%  "
$ws.Range("E23").Value = "  -3.10%  "
$ws.Range("E24").Value = "  -1.50%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.00"
$ws.Range("E26").Value = "  +0.00%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.20"
$ws.Range("E27").Value = "  -2.75%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0₃0859"
$ws.Range("E28").Value = "  -5.88%  "
$ws.Range("E29").Value = "  +1.95%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.24"
$ws.Range("E30").Value = "  +0.10%  "
$ws.Range("E31").Value = "  -1.25%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "166.22"
$ws.Range("E32").Value = "  +1.53%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.48"
$ws.Range("E33").Value = "  -0.02%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.87"
$ws.Range("E34").Value = "  -0.57%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.999"
$ws.Range("E35").Value = "  +0.02%  "
$ws.Range("E36").Value = "  -2.44%  "
$ws.Range("E37").Value = "  -1.26%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "349.81"
$ws.Range("E38").Value = "  -0.32%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.966"
$ws.Range("E39").Value = "  -2.39%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.31"
$ws.Range("E40").Value = "  -0.71%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "4.01"
$ws.Range("E41").Value = "  -2.66%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "38.31"
$ws.Range("E42").Value = "  -0.98%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "20.94"
$ws.Range("E43").Value = "  -4.69%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "20.35"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0563"
$ws.Range("E45").Value = "  -3.60%  "
$ws.Range("E46").Value = "  -1.24%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.999"
$ws.Range("E47").Value = "  +0.06%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "11.02"
$ws.Range("E48").Value = "  -0.06%  "
$ws.Range("E49").Value = "  -3.10%  "
$ws.Range("E50").Value = "  -2.97%  "
$ws.Range("B51").Value = "Aave"
$ws.Range("C51").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "128.49"
$ws.Range("E51").Value = "  -4.55%  "
